$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.846.45"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "3.847.57"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "698.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("D7").Value = "3.846.20"
$ws.Range("E7").Value = "  +1.65%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.60%  "

$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "4.500.13"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").Value = "3.846.17"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("D17").Value = "70.925.09"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.114"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "494.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.23%  "

$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.06%  "

$ws.Range("E26").Value = "  +0.67%  "

$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("E28").Value = "  -2.87%  "

$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E34").Value = "  +2.04%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.66%  "

$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.807.38"
$ws.Range("E36").Value = "  +1.98%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.43%  "

$ws.Range("E40").Value = "  +7.81%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.07%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +1.86%  "

$ws.Range("E46").Value = "  -5.65%  "

$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("E49").Value = "  -3.40%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.17%  "
